# The presentation currently carries the "Integral" (Red Violet) design as
# its main theme (ppt/theme/theme1.xml, used by the Slide Master) while the
# default "Office Theme" colours survive only as the Notes Master's theme
# (ppt/theme/theme2.xml). This edit restores the default "Office Theme"
# palette as the presentation's main design by rewriting the twelve
# theme colours (via the modern ThemeColorScheme, which is what actually
# backs <a:clrScheme> in the OOXML) on the shared Slide Master theme.
#
# NOTE: PowerPoint's ColorFormat/RGB COM properties store colours as
# 0x00BBGGRR (the classic Win32 RGB() macro ordering), so every target
# "RRGGBB" hex value below is byte-swapped before being assigned.

$p = $ppt.ActivePresentation

# Any slide can reach the shared Slide Master theme colour scheme - they
# all point at the same underlying theme part, so one slide is enough.
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Target palette = the stock Office Theme colours (dk1/lt1 = black/white
# are already correct and shared by both palettes, so they're left alone).
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
